$d = $word.ActiveDocument
Write-Output $d.Content.Text.Substring(0, 200)
